$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New bug report (BR-003-ish, second block rows 25-38): words.js second-parameter bug ---

$ws.Range("B25").Value = "Possible programming error in the file words.js, when words() is used with a string parameter. Refers to Test ID TC022"
$ws.Rows(25).RowHeight = 28.8

$ws.Range("B26").Value = "When the words function is given a second parameter as a string, the function returns an unexpected array. "
$ws.Rows(26).RowHeight = 28.8

$ws.Range("B27").Value = "petrikreus"

$ws.Range("B28").Value = 44910

$ws.Range("B29").Value = "COMP.SE.200-2022-2023-1 Utils function library"

$ws.Range("B30").Value = "Node v16.17.0, Jest 29.1.2, Coveralls 3.1.1, Ubuntu 20.04.5 LTS (WSL2)"

$ws.Range("B33").Value = "Minor"

$ws.Range("B34").Value = "Locally: npm test"

$ws.Range("B35").Value = "1. Force the usage with RegExp by throwing an error when a string or number is used as a second parameter.
2. Check if the second parameter is a string. If it is, escape special characters to avoid unexpected results and return first index of the result array: string.match(pattern)[0].
See: https://developer.mozilla.org/en-US/docs/Web/JavaScript/Reference/Global_Objects/String/match#a_non-regexp_as_the_parameter"
$ws.Rows(35).RowHeight = 144

$ws.Range("B36").Value = "otula"

$ws.Range("B38").Value = "The main use case of the function is expected to be without the second parameter, or the second parameter being a RegExp. That is the reason for minor seriousness."
$ws.Rows(38).RowHeight = 43.2

# --- Style tweak: left-align the date-formatted cells (B8 / B28 share this style) ---
$ws.Range("B8").NumberFormat = "m/d/yyyy"
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B28").NumberFormat = "m/d/yyyy"
$ws.Range("B28").HorizontalAlignment = -4131

# --- View state: scrolled down to the new bug report, selection parked past the data ---
$ws.Range("A50").Select()
$excel.ActiveWindow.ScrollRow = 23
